$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.036839623004198074
$ws.Range("C2").Value = 0.016010822728276253
$ws.Range("D2").Value = 0.011369539424777031
$ws.Range("E2").Value = 0.008069007657468319
$ws.Range("F2").Value = 0.00006270490848692134
$ws.Range("J2").Value = 0.12766611576080322
$ws.Range("K2").Value = 1.4641062021255493
